$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.530.21"
$ws.Range("E2").Value = "  +0.62%  "
$ws.Range("D3").Value = "1.729.74"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9993"
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.77"
$ws.Range("E5").Value = "  +1.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4812"
$ws.Range("E7").Value = "  +1.49%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2676"
$ws.Range("E8").Value = "  +1.41%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06173"
$ws.Range("E9").Value = "  -0.46%  "
$ws.Range("D10").Value = "1.728.83"
$ws.Range("E10").Value = "  +0.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07187"
$ws.Range("E11").Value = "  +1.70%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.56"
$ws.Range("E12").Value = "  +0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6093"
$ws.Range("E13").Value = "  +1.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.530"
$ws.Range("E14").Value = "  +2.36%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.23"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "26.528.67"
$ws.Range("E17").Value = "  +0.57%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9996"
$ws.Range("E18").Value = "  +0.04%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006940"
$ws.Range("E19").Value = "  +1.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  -0.23%  "
$ws.Range("D21").Value = "1.951.25"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.527"
$ws.Range("E22").Value = "  -0.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.805"
$ws.Range("E23").Value = "  +0.59%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.251"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "136.99"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.34"
$ws.Range("E26").Value = "  +0.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.781"
$ws.Range("E27").Value = "  +0.59%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.408"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.25"
$ws.Range("E29").Value = "  +0.22%  "
$ws.Range("E30").Value = "  +0.49%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08025"
$ws.Range("E31").Value = "  +2.80%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.699"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04510"
$ws.Range("E33").Value = "  +0.24%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.616"
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.007"
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6259"
$ws.Range("E36").Value = "  +0.39%  "
$ws.Range("B37").Value = "RenderToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.085"
$ws.Range("E37").Value = "  +7.32%  "
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9074"
$ws.Range("E38").Value = "  -3.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.388"
$ws.Range("E39").Value = "  -2.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.003"
$ws.Range("E40").Value = "  +0.24%  "
$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01502"
$ws.Range("E41").Value = "  +0.87%  "
$ws.Range("B42").Value = "Quant"
$ws.Range("C42").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "102.49"
$ws.Range("E42").Value = "  -9.84%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.540"
$ws.Range("E43").Value = "  -2.30%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.966"
$ws.Range("E45").Value = "  +9.48%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1179"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05378"
$ws.Range("E47").Value = "  +1.96%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.74"
$ws.Range("E48").Value = "  +1.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.821"
$ws.Range("E49").Value = "  -0.17%  "
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("E51").Value = "  +0.42%  "
